# DO + EFO feb 2023 update
# Update the "version" column (E) for the Disease Ontology (DO) and
# Experimental Factor Ontology (EFO) rows in the metadata sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Disease Ontology -> version v2023-01-31 becomes v2023-02-27
$ws.Range("E3").Value = "v2023-02-27"

# Row 4: Experimental Factor Ontology -> version v3.50.0 becomes v3.51.0
$ws.Range("E4").Value = "v3.51.0"

# Move the active selection to E3 (matches the saved cursor position)
[void]$ws.Range("E3").Select()
